$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DMD")
$ws.Activate()

# Clear the contents (values/formulas) of the small "TimeSlice" helper table
# while leaving the existing cell formatting/styles untouched.
$ws.Range("B22:D32").ClearContents()
$ws.Range("D34").ClearContents()

# Restore the view state (scroll position / zoom / selection) that Excel
# saved along with the content edit.
$ws.Application.ActiveWindow.Zoom = 240
$ws.Range("A18").Select()
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("G27").Select()
